$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rotate the "Robyn Guarriello/Her Mine E" (row 3), "Disha Jain/Ron Ferretly" (row 5)
# and "Rohni Awasthi/Larry Richards" (row 6) pairs so that:
#   row 3 <- old row 6
#   row 5 <- old row 3
#   row 6 <- old row 5
$row6A = $ws.Cells.Item(6, 1).Value2
$row6B = $ws.Cells.Item(6, 2).Value2
$row3A = $ws.Cells.Item(3, 1).Value2
$row3B = $ws.Cells.Item(3, 2).Value2
$row5A = $ws.Cells.Item(5, 1).Value2
$row5B = $ws.Cells.Item(5, 2).Value2

$ws.Cells.Item(3, 1).Value = $row6A
$ws.Cells.Item(3, 2).Value = $row6B

$ws.Cells.Item(5, 1).Value = $row3A
$ws.Cells.Item(5, 2).Value = $row3B

$ws.Cells.Item(6, 1).Value = $row5A
$ws.Cells.Item(6, 2).Value = $row5B
